$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45967, 8609.04430483594, 8386.819082362,   13084, 4024.3820024425,  -28.0332881331459),
    @(45968, 8502.80245187863, 7509.55313039134, 9004,  4601.89725334072, 129.477099322169),
    @(45969, 2984.69117368278, 5087.72996960059, 9004,  4177.27218806039, 10.8750899025409),
    @(45970, 2763.1996804858,  5039.8342118635,  9004,  4070.2410604268,  4.41980301209613),
    @(45971, 8898.0650431158,  8344.99900957978, 9004,  5000.39359711733, 180.891358612379),
    @(45972, 9273.47299076555, 8934.64107360856, 9004,  5243.45317921512, 215.58726053432),
    @(45973, 9273.47299076555, 8845.11722633657, 9004,  5243.45317921512, 211.857100231321),
    @(45974, 9273.47299076555, 8676.99421568198, 9004,  5243.45317921512, 204.851974787379),
    @(45975, 9273.47299076555, 8203.19159203698, 9004,  5243.45317921512, 185.110198802171),
    @(45976, 3551.20396235217, 5769.88654315521, 9004,  4928.56549759564, 70.6021683646188),
    @(45977, 3412.56435984501, 5437.57477329455, 9004,  4921.09368438029, 56.4445190697847),
    @(45978, 10116.508538002,  8792.85308422125, 9004,  5829.29824723688, 234.089638810756),
    @(45979, 10116.508538002,  9018.34316514545, 9004,  5829.29824723688, 243.485058849264),
    @(45980, 10116.508538002,  8923.74917223295, 9004,  5829.29824723688, 239.543642477909)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
